$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price values so they stay as text (matches source inlineStr cells)
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply new values
$ws.Range("D2").Value = '26.385.72'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '1.695.25'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").Value = '218.18'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").Value = '0.5429'
$ws.Range("E6").Value = '  +3.23%  '
$ws.Range("D7").Value = '1.010'
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("D8").Value = '0.2733'
$ws.Range("E8").Value = '  +0.59%  '
$ws.Range("D9").Value = '0.06448'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").Value = '21.89'
$ws.Range("E10").Value = '  -0.96%  '
$ws.Range("D11").Value = '0.07674'
$ws.Range("E11").Value = '  +2.16%  '
$ws.Range("D12").Value = '1.696.94'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D13").Value = '4.551'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = '0.5834'
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("D15").Value = '0.000008354'
$ws.Range("E15").Value = '  -2.01%  '
$ws.Range("D16").Value = '66.07'
$ws.Range("E16").Value = '  +2.25%  '
$ws.Range("D17").Value = '26.442.98'
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").Value = '4.939'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").Value = '10.93'
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("D21").Value = '190.80'
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("D22").Value = '6.266'
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("D23").Value = '1.012'
$ws.Range("E23").Value = '  +0.49%  '
$ws.Range("D24").Value = '148.51'
$ws.Range("D25").Value = '0.1305'
$ws.Range("E25").Value = '  +4.68%  '
$ws.Range("D26").Value = '7.919'
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("D27").Value = '15.79'
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '1.392'
$ws.Range("E28").Value = '  +2.55%  '
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value = '0.06216'
$ws.Range("E29").Value = '  -6.79%  '
$ws.Range("D30").Value = '1.327'
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("D31").Value = '3.610'
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").Value = '3.586'
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("D33").Value = '1.698'
$ws.Range("E33").Value = '  +2.12%  '
$ws.Range("D34").Value = '1.037'
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("D35").Value = '0.6149'
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("D36").Value = '2.412'
$ws.Range("E36").Value = '  +0.50%  '
$ws.Range("D37").Value = '2.760'
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("D38").Value = '0.01648'
$ws.Range("E38").Value = '  +1.51%  '
$ws.Range("D39").Value = '1.115.21'
$ws.Range("E39").Value = '  +0.52%  '
$ws.Range("D40").Value = '6.101'
$ws.Range("E40").Value = '  -5.60%  '
$ws.Range("D41").Value = '0.8846'
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").Value = '1.015'
$ws.Range("D43").Value = '101.31'
$ws.Range("E43").Value = '  +0.52%  '
$ws.Range("D44").Value = '1.848.52'
$ws.Range("E44").Value = '  +0.82%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '57.74'
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("D47").Value = '8.203'
$ws.Range("E47").Value = '  +0.65%  '
$ws.Range("D48").Value = '1.006'
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("D49").Value = '0.05294'
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("D50").Value = '6.114'
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("D51").Value = '0.4302'
$ws.Range("E51").Value = '  +0.11%  '
